$d = $word.ActiveDocument

# --- Simple single-value cell replacements (unique text, safe with Find/Execute) ---
$d.Content.Find.Execute("93.48", $true, $false, $false, $false, $false, $true, 1, $false, "0M", 2)
$d.Content.Find.Execute("116.22", $true, $false, $false, $false, $false, $true, 1, $false, "0M", 2)
$d.Content.Find.Execute("1782", $true, $false, $false, $false, $false, $true, 1, $false, "0M", 2)
$d.Content.Find.Execute("33066", $true, $false, $false, $false, $false, $true, 1, $false, "36270", 2)
$d.Content.Find.Execute("0.00633", $true, $false, $false, $false, $false, $true, 1, $false, "0.01334", 2)
$d.Content.Find.Execute("0.00555", $true, $false, $false, $false, $false, $true, 1, $false, "0.00583", 2)
$d.Content.Find.Execute("0.01397", $true, $false, $false, $false, $false, $true, 1, $false, "0.03197", 2)
$d.Content.Find.Execute("0.02237", $true, $false, $false, $false, $false, $true, 1, $false, "0.03306", 2)
$d.Content.Find.Execute("11.66797", $true, $false, $false, $false, $false, $true, 1, $false, "116.06418", 2)

# --- Collapse the three multi-run "raw samples" cells near the end of the table
#     down to a single summary value (removes the extra runs/tabs) ---
$tbl = $d.Tables.Item(1)
$tbl.Cell(44, 1).Range.Text = "93.48"
$tbl.Cell(45, 1).Range.Text = "116.22"
$tbl.Cell(46, 1).Range.Text = "1782"
